# Updates cryptos list price/volume(1h) data (and row 45/46 coin swap)
# to match the Sun Apr 16 03:41:39 UTC 2023 GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.553.44'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.109.11'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.011'
$ws.Range('E4').Value = '  +0.75%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '335.84'
$ws.Range('E5').Value = '  +1.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.009'
$ws.Range('E6').Value = '  +0.69%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5243'
$ws.Range('E7').Value = '  +0.58%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4545'
$ws.Range('E8').Value = '  +4.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '55.48'
$ws.Range('E9').Value = '  +2.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09023'
$ws.Range('E10').Value = '  +1.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.169'
$ws.Range('E11').Value = '  +1.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.59'
$ws.Range('E12').Value = '  +0.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.112.73'
$ws.Range('E13').Value = '  +1.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.845'
$ws.Range('E14').Value = '  +2.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.111'
$ws.Range('E15').Value = '  +5.78%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001176'
$ws.Range('E16').Value = '  +4.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '97.08'
$ws.Range('E17').Value = '  +1.30%  '
$ws.Range('E18').Value = '  +0.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06683'
$ws.Range('E19').Value = '  +1.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.33'
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('E21').Value = '  +0.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.245'
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.616.69'
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.79'
$ws.Range('E24').Value = '  +4.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.358'
$ws.Range('E25').Value = '  +0.69%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.367.22'
$ws.Range('E26').Value = '  +1.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.35'
$ws.Range('E27').Value = '  +0.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '163.29'
$ws.Range('E28').Value = '  +0.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.510'
$ws.Range('E29').Value = '  -2.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.56'
$ws.Range('E30').Value = '  +1.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.212'
$ws.Range('E31').Value = '  +2.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.343'
$ws.Range('E33').Value = '  +3.37%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.619'
$ws.Range('E34').Value = '  -1.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.962'
$ws.Range('E35').Value = '  +1.74%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.40'
$ws.Range('E36').Value = '  +3.92%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.879'
$ws.Range('E37').Value = '  +8.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02614'
$ws.Range('E38').Value = '  +1.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06805'
$ws.Range('E39').Value = '  -0.06%  '
$ws.Range('E40').Value = '  +2.39%  '
$ws.Range('E41').Value = '  -1.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6835'
$ws.Range('E42').Value = '  -0.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.259'
$ws.Range('E43').Value = '  -0.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6434'
$ws.Range('E44').Value = '  +1.11%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.05'
$ws.Range('E45').Value = '  +0.73%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.303'
$ws.Range('E46').Value = '  +4.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.672'
$ws.Range('E47').Value = '  +1.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.00000000354'
$ws.Range('E48').Value = '  +18.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.250'
$ws.Range('E49').Value = '  +0.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.209'
$ws.Range('E50').Value = '  -2.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '83.18'
$ws.Range('E51').Value = '  +1.79%  '
